$d = $word.ActiveDocument

function Find-ParagraphByText($doc, $text) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.Trim() -eq $text) {
            return $p
        }
    }
    return $null
}

$pkgHeader = "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:body>"
$pkgFooter = "</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"

# ---------------------------------------------------------------------------
# 1) "Non-Functional Requirements" paragraph: add <w:lastRenderedPageBreak/>
#    before the text run.
# ---------------------------------------------------------------------------
$p1 = Find-ParagraphByText $d "Non-Functional Requirements"
$r1 = $p1.Range
$r1.Collapse(1)
$xml1 = $pkgHeader + "<w:p><w:pPr><w:pStyle w:val='Heading2'/><w:numPr><w:ilvl w:val='1'/><w:numId w:val='9'/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Non-Functional Requirements</w:t></w:r></w:p>" + $pkgFooter
$r1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# 2) The "Object O" / "riented Analysis (OOA) - UML" split-run heading
#    becomes a single run, and the _GoBack bookmark that used to sit
#    between those two runs is removed from here (it gets re-added at the
#    end of the "Analysis Classes" bullet below).
# ---------------------------------------------------------------------------
$p2 = Find-ParagraphByText $d "Object Oriented Analysis (OOA) – UML"
$em = [string][char]0x2013
$p2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Object O")) {
        $p2 = $p
        break
    }
}
$r2 = $p2.Range
$r2.Collapse(1)
$xml2 = $pkgHeader + "<w:p><w:pPr><w:pStyle w:val='Heading1'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='9'/></w:numPr></w:pPr><w:r><w:t>Object Oriented Analysis (OOA) " + $em + " UML</w:t></w:r></w:p>" + $pkgFooter
$r2.InsertXML($xml2)

# ---------------------------------------------------------------------------
# 3) Insert a new "Use Case Expansion" Heading3 bullet right after
#    "Use Cases" (this is the new 4th bullet under section 4 / OOA).
# ---------------------------------------------------------------------------
$idx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.Trim() -eq "Use Cases") {
        $idx = $i
        break
    }
}
$p3 = $d.Paragraphs($idx)
$p3.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($idx + 1)
$newPara.Range.InsertAfter("Use Case Expansion")

# ---------------------------------------------------------------------------
# 4) Re-add the _GoBack bookmark around the end of "Analysis Classes".
# ---------------------------------------------------------------------------
$p4 = Find-ParagraphByText $d "Analysis Classes"
$r4 = $p4.Range
$r4.Collapse(1)
$xml4 = $pkgHeader + "<w:p><w:pPr><w:pStyle w:val='Heading3'/><w:numPr><w:ilvl w:val='2'/><w:numId w:val='9'/></w:numPr></w:pPr><w:r><w:t>Analysis Classes</w:t></w:r><w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/></w:p>" + $pkgFooter
$r4.InsertXML($xml4)

Write-Host "done"
